$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '协鑫集成'
$ws.Cells.Item(2, 2).Value = '协鑫集成'
$ws.Cells.Item(2, 3).Value = '杉杉股份'

$ws.Cells.Item(3, 1).Value = '数据港'
$ws.Cells.Item(3, 2).Value = '杉杉股份'
$ws.Cells.Item(3, 3).Value = '协鑫集成'

$ws.Cells.Item(4, 1).Value = '神剑股份'
$ws.Cells.Item(4, 2).Value = '数据港'
$ws.Cells.Item(4, 3).Value = '永太科技'

$ws.Cells.Item(5, 1).Value = '杭电股份'
$ws.Cells.Item(5, 2).Value = '神剑股份'
$ws.Cells.Item(5, 3).Value = '神剑股份'

$ws.Cells.Item(6, 1).Value = '杉杉股份'
$ws.Cells.Item(6, 2).Value = '杭电股份'
$ws.Cells.Item(6, 3).Value = '巨力索具'

$ws.Cells.Item(7, 1).Value = '五洲新春'
$ws.Cells.Item(7, 2).Value = '巨力索具'
$ws.Cells.Item(7, 3).Value = '杭电股份'

$ws.Cells.Item(8, 1).Value = '巨力索具'
$ws.Cells.Item(8, 2).Value = '银河电子'
$ws.Cells.Item(8, 3).Value = '数据港'

$ws.Cells.Item(9, 1).Value = '航天发展'
$ws.Cells.Item(9, 2).Value = '协鑫能科'
$ws.Cells.Item(9, 3).Value = '利欧股份'

$ws.Cells.Item(10, 1).Value = '银河电子'
$ws.Cells.Item(10, 2).Value = '五洲新春'
$ws.Cells.Item(10, 3).Value = '浙文互联'

$ws.Cells.Item(11, 1).Value = '利欧股份'
$ws.Cells.Item(11, 2).Value = '洲际油气'
$ws.Cells.Item(11, 3).Value = '雷科防务'

$ws.Cells.Item(12, 1).Value = '永太科技'
$ws.Cells.Item(12, 2).Value = '天奇股份'
$ws.Cells.Item(12, 3).Value = '银河电子'

$ws.Cells.Item(13, 1).Value = '洲际油气'
$ws.Cells.Item(13, 2).Value = '三变科技'
$ws.Cells.Item(13, 3).Value = '航天发展'

$ws.Cells.Item(14, 1).Value = '雷科防务'
$ws.Cells.Item(14, 2).Value = '利欧股份'
$ws.Cells.Item(14, 3).Value = '洲际油气'

$ws.Cells.Item(15, 1).Value = '协鑫能科'
$ws.Cells.Item(15, 2).Value = '雷科防务'
$ws.Cells.Item(15, 3).Value = '天奇股份'

$ws.Cells.Item(16, 1).Value = '天奇股份'
$ws.Cells.Item(16, 2).Value = '永太科技'
$ws.Cells.Item(16, 3).Value = '五洲新春'

$ws.Cells.Item(17, 1).Value = '浙文互联'
$ws.Cells.Item(17, 2).Value = '航天电子'
$ws.Cells.Item(17, 3).Value = '中超控股'

$ws.Cells.Item(18, 1).Value = '三变科技'
$ws.Cells.Item(18, 2).Value = '贵州茅台'
$ws.Cells.Item(18, 3).Value = 'TCL中环'

$ws.Cells.Item(19, 1).Value = '工业富联'
$ws.Cells.Item(19, 2).Value = '浙文互联'
$ws.Cells.Item(19, 3).Value = '白银有色'

$ws.Cells.Item(20, 1).Value = '网宿科技'
$ws.Cells.Item(20, 2).Value = '东方财富'
$ws.Cells.Item(20, 3).Value = '湖南白银'

$ws.Cells.Item(21, 1).Value = '信维通信'
$ws.Cells.Item(21, 2).Value = '百川股份'
$ws.Cells.Item(21, 3).Value = '平潭发展'

